# Sprint 3 Review and burnup, + final documentation
#
# 1) Fill in the "Total Done" (column C) burnup values for Sprint 3 -
#    these were all placeholder 0s and now carry the real cumulative
#    totals, matching the formatting already used by the neighboring
#    "Estimated Project Unit" / "Ideal Done" columns (B & D).
# 2) Apply the final print/page-layout documentation pass to the sheet
#    (portrait orientation, zero header/footer margins).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Burnup ("Total Done") values for C2:C15 -------------------------
$ws.Range("C2").Value  = 1
$ws.Range("C3").Value  = 3
$ws.Range("C4").Value  = 5
$ws.Range("C5").Value  = 8
$ws.Range("C6").Value  = 12
$ws.Range("C7").Value  = 15
$ws.Range("C8").Value  = 18
$ws.Range("C9").Value  = 21
$ws.Range("C10").Value = 24
$ws.Range("C11").Value = 27
$ws.Range("C12").Value = 29
$ws.Range("C13").Value = 30
$ws.Range("C14").Value = 35
$ws.Range("C15").Value = 38

# C14:C15 already carried the same cell format as the B/D columns;
# copy that formatting down over C2:C13 so the whole column is
# consistent (same style index the diff shows for the whole range).
$ws.Range("C14").Copy()
$ws.Range("C2:C13").PasteSpecial(-4122)  # xlPasteFormats

# --- Final page setup / print documentation ---------------------------
$ws.PageSetup.Orientation   = 1  # xlPortrait
$ws.PageSetup.HeaderMargin  = 0
$ws.PageSetup.FooterMargin  = 0
